$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.499.33'
$ws.Range('E2').Value = '  +3.13%  '
$ws.Range('D3').Value = '3.070.08'
$ws.Range('E3').Value = '  +2.22%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '551.57'
$ws.Range('E5').Value = '  +2.59%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.44'
$ws.Range('E6').Value = '  +5.78%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '3.066.39'
$ws.Range('E8').Value = '  +2.28%  '
$ws.Range('E9').Value = '  +1.28%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.51'
$ws.Range('E10').Value = '  +5.90%  '
$ws.Range('E11').Value = '  +2.40%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.458'
$ws.Range('E12').Value = '  +2.24%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000228'
$ws.Range('E13').Value = '  +2.65%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.97'
$ws.Range('E14').Value = '  +2.89%  '
$ws.Range('D15').Value = '3.566.59'
$ws.Range('E15').Value = '  +2.35%  '
$ws.Range('D16').Value = '63.480.08'
$ws.Range('E16').Value = '  +3.14%  '
$ws.Range('D17').Value = '3.073.22'
$ws.Range('E17').Value = '  +2.37%  '
$ws.Range('E18').Value = '  -1.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.80'
$ws.Range('E19').Value = '  +2.49%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '485.69'
$ws.Range('E20').Value = '  +3.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.93'
$ws.Range('E21').Value = '  +5.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.679'
$ws.Range('E22').Value = '  +0.44%  '
$ws.Range('E23').Value = '  +5.25%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '81.19'
$ws.Range('E24').Value = '  +1.09%  '
$ws.Range('E25').Value = '  +6.58%  '
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('E27').Value = '  +3.50%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.90'
$ws.Range('E28').Value = '  +1.48%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.02'
$ws.Range('E29').Value = '  +7.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  -0.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '26.26'
$ws.Range('E31').Value = '  +2.36%  '
$ws.Range('E32').Value = '  +1.29%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.45'
$ws.Range('E33').Value = '  +7.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.72'
$ws.Range('E34').Value = '  +3.89%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '55.56'
$ws.Range('E35').Value = '  +1.22%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.01'
$ws.Range('E36').Value = '  +1.67%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '467.76'
$ws.Range('E37').Value = '  +2.83%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0825'
$ws.Range('E38').Value = '  +4.76%  '
$ws.Range('E39').Value = '  +3.79%  '
$ws.Range('D40').Value = '3.041.92'
$ws.Range('E40').Value = '  -4.09%  '
$ws.Range('E41').Value = '  -1.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.24'
$ws.Range('E42').Value = '  +1.31%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.58'
$ws.Range('E43').Value = '  +4.56%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '27.74'
$ws.Range('E44').Value = '  +2.74%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.256'
$ws.Range('E45').Value = '  +4.77%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.05'
$ws.Range('E47').Value = '  +2.39%  '
$ws.Range('E48').Value = '  +2.33%  '
$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '116.94'
$ws.Range('E49').Value = '  -1.83%  '
$ws.Range('B50').Value = 'PEPE'
$ws.Range('C50').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D50').Value = '0.0₃0511'
$ws.Range('E50').Value = '  +3.04%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.09'
$ws.Range('E51').Value = '  +4.29%  '
